# fixing author calc of total churn logic in authors_contrib function
#
# The "file_legacy_complexity" sheet previously had only one row per file
# (using whichever author/legacy-percentage happened to be last computed).
# It should instead have one row per (file, author) pair, mirroring the
# "file_author_contrib" sheet, with the correct per-author legacy
# percentage (column D) and a recomputed "cog_complexity" churn value
# (column E).

$wb = $excel.ActiveWorkbook
$ws3 = $wb.Sheets.Item("file_author_contrib")
$ws4 = $wb.Sheets.Item("file_legacy_complexity")

# file_name (B) and author (C) line up 1-for-1 with file_author_contrib's
# rows 2:20 -- reuse them directly (via copy/paste) so the shared-string
# references are preserved exactly instead of being re-typed.
$ws3.Range("B2:C20").Copy()
$ws4.Range("B2").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# id (A)
for ($r = 2; $r -le 20; $r++) {
    $ws4.Cells.Item($r, 1).Value = $r - 1
}

# legacy_percentage (D) -- same values as the "percentages" column (H) on
# file_author_contrib.
$dValues = @(
    1, 1, 1, 1, 1,
    0.8666666666666667, 0.1333333333333333,
    0.5933503836317136, 0.4066496163682864,
    1, 1, 1, 1, 1, 1,
    0.9415584415584416, 0.05844155844155844,
    1, 1
)

# cog_complexity (E) -- total churn attributable to this author on this
# file.
$eValues = @(
    0, 0, 7.333333333333333, 0, 0,
    0, 0,
    6.615384615384615, 6.615384615384615,
    2.5, 0, 4.5, 6, 4, 9.384615384615385,
    5.25, 5.25,
    2.25, 46
)

for ($i = 0; $i -lt 19; $i++) {
    $r = $i + 2
    $ws4.Cells.Item($r, 4).Value = $dValues[$i]
    $ws4.Cells.Item($r, 5).Value = $eValues[$i]
}
